$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as TEXT (keeps numeric-looking strings as strings,
# without leaving a residual number-format style applied to the cell, matching
# the plain <c t="s"> cells used for this data in the workbook).
function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Helper: set a cell value as a NUMBER
function Set-NumCell($addr, $val) {
    $ws.Range($addr).Value = $val
}

# The order in which brand-new text values are first written controls the
# order they are appended to the shared-strings table, so cells are touched
# in a specific sequence below (matching how the workbook was authored)
# rather than in simple row/column order.

Set-TextCell "F4" "K"
Set-TextCell "G4" "632"
Set-TextCell "B2" "2200"
Set-TextCell "C2" "SE"
Set-TextCell "D2" "8123730680"
Set-TextCell "H2" "00296"
Set-TextCell "I2" "0000385977"
Set-TextCell "K2" "Charlie"
Set-TextCell "L2" "Charlie 22"
Set-TextCell "M2" "I638176"
Set-TextCell "N2" "2305PXT6252"
Set-TextCell "E3" "9000667710"
Set-TextCell "D3" "8123731130"
Set-TextCell "G3" "161"
Set-TextCell "I3" "0000007905"

Set-TextCell "G2" "632"
Set-TextCell "B3" "2200"
Set-TextCell "C3" "SE"
Set-TextCell "H3" "00296"
Set-TextCell "J3" "FERT"
Set-TextCell "K3" "Charlie"
Set-TextCell "L3" "Charlie 22"
Set-TextCell "M3" "I638176"
Set-TextCell "N3" "2305PXT6252"
Set-TextCell "B4" "2200"
Set-TextCell "C4" "SE"
Set-TextCell "D4" "8123730680"
Set-TextCell "E4" "-2"
Set-TextCell "I4" "0000385977"
Set-TextCell "J4" "FERT"
Set-TextCell "K4" "Charlie"
Set-TextCell "L4" "Charlie 22"
Set-TextCell "M4" "I638176"
Set-TextCell "N4" "2305PXT6252"

# ---- numeric cells ----
Set-NumCell "A2" 44950
Set-NumCell "O2" 1
Set-NumCell "P2" 769.12

Set-NumCell "A3" 44950
Set-NumCell "O3" -1
Set-NumCell "P3" 769.12

Set-NumCell "A4" 44950
Set-NumCell "O4" -1
Set-NumCell "P4" 769.12

# A4 gets the same date/time display style as A2 / A3
$ws.Range("A4").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

# ---- sheet-level formatting tweaks ----
$ws.Range("Q14").Select() | Out-Null

# Column width for column A (closest reachable value to the authored width)
$ws.Columns.Item(1).ColumnWidth = 17.25
